$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.983.38"
$ws.Range("E2").Value = "'  -1.58%  "

$ws.Range("D3").Value = "'1.974.46"
$ws.Range("E3").Value = "'  -3.38%  "

$ws.Range("E4").Value = "'  -0.01%  "

$ws.Range("D5").Value = "'239.01"
$ws.Range("E5").Value = "'  -7.17%  "

$ws.Range("D6").Value = "'0.597"
$ws.Range("E6").Value = "'  -4.47%  "

$ws.Range("E7").Value = "'  +0.05%  "

$ws.Range("D8").Value = "'53.82"
$ws.Range("E8").Value = "'  -7.11%  "

$ws.Range("B9").Value = "'OKB"
$ws.Range("C9").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'59.15"
$ws.Range("E9").Value = "'  +3.37%  "

$ws.Range("B10").Value = "'Cardano"
$ws.Range("C10").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.368"
$ws.Range("E10").Value = "'  -5.51%  "

$ws.Range("B11").Value = "'Dogecoin"
$ws.Range("C11").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0744"
$ws.Range("E11").Value = "'  -7.77%  "

$ws.Range("B12").Value = "'TRON"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.0983"
$ws.Range("E12").Value = "'  -5.05%  "

$ws.Range("B13").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "'2.261.35"
$ws.Range("E13").Value = "'  -3.54%  "

$ws.Range("B14").Value = "'Chainlink"
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'13.82"
$ws.Range("E14").Value = "'  -7.43%  "

$ws.Range("B15").Value = "'Avalanche"
$ws.Range("C15").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'20.81"
$ws.Range("E15").Value = "'  -3.44%  "

$ws.Range("B16").Value = "'Polygon"
$ws.Range("C16").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.746"
$ws.Range("E16").Value = "'  -9.74%  "

$ws.Range("B17").Value = "'Polkadot"
$ws.Range("C17").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'5.00"
$ws.Range("E17").Value = "'  -7.86%  "

$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'1.974.39"
$ws.Range("E18").Value = "'  -3.49%  "

$ws.Range("B19").Value = "'WrappedBTC"
$ws.Range("C19").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "'36.871.95"
$ws.Range("E19").Value = "'  -1.53%  "

$ws.Range("B20").Value = "'Litecoin"
$ws.Range("C20").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'67.70"
$ws.Range("E20").Value = "'  -3.83%  "

$ws.Range("B21").Value = "'ShibaInu"
$ws.Range("C21").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0₃0803"
$ws.Range("E21").Value = "'  -6.79%  "

$ws.Range("B22").Value = "'BitcoinCash"
$ws.Range("C22").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'226.37"
$ws.Range("E22").Value = "'  -1.54%  "

$ws.Range("B23").Value = "'Uniswap"
$ws.Range("C23").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.91"
$ws.Range("E23").Value = "'  -7.21%  "

$ws.Range("B24").Value = "'Dai"
$ws.Range("C24").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "'  -0.04%  "

$ws.Range("E25").Value = "'  -1.55%  "

$ws.Range("B26").Value = "'PancakeSwap"
$ws.Range("C26").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "'  -13.94%  "

$ws.Range("B27").Value = "'Monero"
$ws.Range("C27").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'161.44"
$ws.Range("E27").Value = "'  -1.59%  "

$ws.Range("B28").Value = "'Cosmos"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'8.53"
$ws.Range("E28").Value = "'  -7.78%  "

$ws.Range("B29").Value = "'EthereumClassic"
$ws.Range("C29").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'18.96"
$ws.Range("E29").Value = "'  -5.52%  "

$ws.Range("B30").Value = "'Kaspa"
$ws.Range("C30").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.126"
$ws.Range("E30").Value = "'  -10.91%  "

$ws.Range("B31").Value = "'ImmutableX"
$ws.Range("C31").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.27"
$ws.Range("E31").Value = "'  -6.69%  "

$ws.Range("B32").Value = "'Stellar"
$ws.Range("C32").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.116"
$ws.Range("E32").Value = "'  -4.04%  "

$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.36"
$ws.Range("E33").Value = "'  -9.64%  "

$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0609"
$ws.Range("E34").Value = "'  -9.15%  "

$ws.Range("B35").Value = "'InternetComputer(DFINITY)"
$ws.Range("C35").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.23"
$ws.Range("E35").Value = "'  -6.61%  "

$ws.Range("B36").Value = "'LidoDAOToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.30"
$ws.Range("E36").Value = "'  -8.12%  "

$ws.Range("B37").Value = "'BinanceUSD"
$ws.Range("C37").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "'  -0.08%  "

$ws.Range("B38").Value = "'WEMIXToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "'  -2.01%  "

$ws.Range("B39").Value = "'RenderToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'3.31"
$ws.Range("E39").Value = "'  -5.72%  "

$ws.Range("B40").Value = "'THORChain"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "'5.08"
$ws.Range("E40").Value = "'  -6.31%  "

$ws.Range("B41").Value = "'HuobiToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'3.02"
$ws.Range("E41").Value = "'  -0.49%  "

$ws.Range("B42").Value = "'Maker"
$ws.Range("C42").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.403.19"
$ws.Range("E42").Value = "'  -0.24%  "

$ws.Range("B43").Value = "'TrustWalletToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.12"
$ws.Range("E43").Value = "'  -7.93%  "

$ws.Range("B44").Value = "'Cronos"
$ws.Range("C44").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0886"
$ws.Range("E44").Value = "'  -9.01%  "

$ws.Range("B45").Value = "'VeChain"
$ws.Range("C45").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0202"
$ws.Range("E45").Value = "'  -7.73%  "

$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'86.48"
$ws.Range("E46").Value = "'  -6.26%  "

$ws.Range("B47").Value = "'InjectiveProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'15.13"
$ws.Range("E47").Value = "'  -8.23%  "

$ws.Range("B48").Value = "'ARBITRUM"
$ws.Range("C48").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.989"
$ws.Range("E48").Value = "'  -6.52%  "

$ws.Range("B49").Value = "'MXToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.85"
$ws.Range("E49").Value = "'  -1.52%  "

$ws.Range("B50").Value = "'FTXToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "'3.62"
$ws.Range("E50").Value = "'  +10.65%  "

$ws.Range("B51").Value = "'FraxShare"
$ws.Range("C51").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'6.57"
$ws.Range("E51").Value = "'  -12.12%  "
